$wb = $excel.ActiveWorkbook

# --- Sheet "Results": add row 6 (written first so shared strings line up
#     with the order Excel produced them in) ---
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("A6").Value = "ship-question-mode-a2-200-log csv2"
$wsResults.Range("B6").Value = "_trimmed-2.csv"
$wsResults.Range("C6").Value = 28.549600000000002
$wsResults.Range("D6").Value = 2.8549600000000002
$wsResults.Range("E6").Value = 5

# match formatting used by the row above it (row 5) -- these columns are
# all plain "centered, General format" cells (style index 1)
$wsResults.Range("A6").HorizontalAlignment = $wsResults.Range("A5").HorizontalAlignment
$wsResults.Range("B6").HorizontalAlignment = $wsResults.Range("B5").HorizontalAlignment
$wsResults.Range("C6").HorizontalAlignment = $wsResults.Range("C5").HorizontalAlignment
$wsResults.Range("D6").HorizontalAlignment = $wsResults.Range("D5").HorizontalAlignment
$wsResults.Range("E6").HorizontalAlignment = $wsResults.Range("E5").HorizontalAlignment

# --- Sheet "Sizes": add row 8 ---
$wsSizes = $wb.Worksheets.Item("Sizes")
$wsSizes.Range("A8").Value = "message-log-200_trimmed-2.csv"
$wsSizes.Range("B8").Value = 403
$wsSizes.Range("C8").Value = 47122
$wsSizes.Range("D8").Value = "Removed unneeded columns, port-id->next-port"

# match formatting used by the row above it (row 7). Column C (Bytes) uses
# a comma-separated-thousands number format; the rest are plain.
$wsSizes.Range("A8").HorizontalAlignment = $wsSizes.Range("A7").HorizontalAlignment
$wsSizes.Range("B8").HorizontalAlignment = $wsSizes.Range("B7").HorizontalAlignment
$wsSizes.Range("C8").HorizontalAlignment = $wsSizes.Range("C7").HorizontalAlignment
$wsSizes.Range("C8").NumberFormat = $wsSizes.Range("C7").NumberFormat
$wsSizes.Range("D8").HorizontalAlignment = $wsSizes.Range("D7").HorizontalAlignment

# --- Update the remembered selection on each sheet ---
$wsSizes.Activate()
$wsSizes.Range("E25").Select()

$wsResults.Activate()
$wsResults.Range("C6").Select()
